$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '43.616.02'
Set-TextValue $ws.Range('E2') '  -1.34%  '
Set-TextValue $ws.Range('D3') '2.233.71'
Set-TextValue $ws.Range('E4') '  +0.14%  '
Set-TextValue $ws.Range('D5') '268.88'
Set-TextValue $ws.Range('E5') '  +3.25%  '
Set-TextValue $ws.Range('D6') '93.10'
Set-TextValue $ws.Range('E6') '  +11.50%  '
Set-TextValue $ws.Range('D7') '0.623'
Set-TextValue $ws.Range('E7') '  -0.96%  '
Set-TextValue $ws.Range('E8') '  +0.05%  '
Set-TextValue $ws.Range('E9') '  +1.79%  '
Set-TextValue $ws.Range('D10') '47.16'
Set-TextValue $ws.Range('E10') '  +5.97%  '
Set-TextValue $ws.Range('D11') '0.0923'
Set-TextValue $ws.Range('E11') '  -1.22%  '
Set-TextValue $ws.Range('D12') '8.33'
Set-TextValue $ws.Range('E12') '  +17.43%  '
Set-TextValue $ws.Range('E13') '  +1.63%  '
Set-TextValue $ws.Range('B14') 'WrappedliquidstakedEther2.0'
Set-TextValue $ws.Range('C14') 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
Set-TextValue $ws.Range('D14') '2.569.23'
Set-TextValue $ws.Range('E14') '  +0.29%  '
Set-TextValue $ws.Range('B15') 'Chainlink'
Set-TextValue $ws.Range('C15') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D15') '15.10'
Set-TextValue $ws.Range('E15') '  +2.74%  '
Set-TextValue $ws.Range('D16') '2.242.66'
Set-TextValue $ws.Range('E16') '  +1.54%  '
Set-TextValue $ws.Range('D17') '0.800'
Set-TextValue $ws.Range('E17') '  +1.90%  '
Set-TextValue $ws.Range('D18') '43.589.02'
Set-TextValue $ws.Range('E18') '  -1.12%  '
Set-TextValue $ws.Range('E19') '  -1.08%  '
Set-TextValue $ws.Range('D20') '6.02'
Set-TextValue $ws.Range('E20') '  -0.38%  '
Set-TextValue $ws.Range('D21') '70.46'
Set-TextValue $ws.Range('E21') '  -1.89%  '
Set-TextValue $ws.Range('E22') '  -1.76%  '
Set-TextValue $ws.Range('D23') '233.23'
Set-TextValue $ws.Range('E23') '  -0.23%  '
Set-TextValue $ws.Range('D24') '8.99'
Set-TextValue $ws.Range('E24') '  -2.82%  '
Set-TextValue $ws.Range('E25') '  +0.00%  '
Set-TextValue $ws.Range('B26') 'Cosmos'
Set-TextValue $ws.Range('C26') 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range('D26') '11.29'
Set-TextValue $ws.Range('E26') '  +4.41%  '
Set-TextValue $ws.Range('B27') 'PancakeSwap'
Set-TextValue $ws.Range('C27') 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Range('D27') '2.49'
Set-TextValue $ws.Range('E27') '  +10.44%  '
Set-TextValue $ws.Range('E28') '  +5.21%  '
Set-TextValue $ws.Range('D29') '40.32'
Set-TextValue $ws.Range('E29') '  -1.26%  '
Set-TextValue $ws.Range('E30') '  +1.91%  '
Set-TextValue $ws.Range('D31') '173.00'
Set-TextValue $ws.Range('E31') '  -0.12%  '
Set-TextValue $ws.Range('D32') '0.0927'
Set-TextValue $ws.Range('E32') '  +3.55%  '
Set-TextValue $ws.Range('D33') '20.83'
Set-TextValue $ws.Range('E33') '  +0.53%  '
Set-TextValue $ws.Range('D34') '5.48'
Set-TextValue $ws.Range('E34') '  +2.07%  '
Set-TextValue $ws.Range('E35') '  +0.27%  '
Set-TextValue $ws.Range('E36') '  -5.24%  '
Set-TextValue $ws.Range('E37') '  -4.51%  '
Set-TextValue $ws.Range('D38') '4.32'
Set-TextValue $ws.Range('E38') '  -4.48%  '
Set-TextValue $ws.Range('D39') '3.60'
Set-TextValue $ws.Range('E39') '  +20.45%  '
Set-TextValue $ws.Range('D40') '12.56'
Set-TextValue $ws.Range('E40') '  -6.51%  '
Set-TextValue $ws.Range('D41') '2.18'
Set-TextValue $ws.Range('E41') '  +2.33%  '
Set-TextValue $ws.Range('D42') '0.220'
Set-TextValue $ws.Range('E42') '  +8.48%  '
Set-TextValue $ws.Range('D43') '63.15'
Set-TextValue $ws.Range('E43') '  -1.36%  '
Set-TextValue $ws.Range('D44') '5.33'
Set-TextValue $ws.Range('E44') '  -4.03%  '
Set-TextValue $ws.Range('D45') '0.0988'
Set-TextValue $ws.Range('E45') '  -0.01%  '
Set-TextValue $ws.Range('D46') '8.39'
Set-TextValue $ws.Range('E46') '  +0.06%  '
Set-TextValue $ws.Range('D47') '100.60'
Set-TextValue $ws.Range('E47') '  -2.72%  '
Set-TextValue $ws.Range('D48') '1.16'
Set-TextValue $ws.Range('E48') '  +2.28%  '
Set-TextValue $ws.Range('E49') '  +2.31%  '
Set-TextValue $ws.Range('D50') '0.436'
Set-TextValue $ws.Range('E50') '  -2.15%  '
Set-TextValue $ws.Range('D51') '2.455.46'
Set-TextValue $ws.Range('E51') '  +0.39%  '
